# Renumber the figure captions on slides 4-6 ("Fig 2/3/4" -> "Fig 3/4/5"),
# splitting the last caption into two runs so the new "Fig 5. " prefix is its
# own run (matching the source edit).

function Get-RestoreEpsilon($originalPts) {
    # Shape geometry (Left/Top/Width/Height) round-trips through a 32-bit
    # "Single" points value. Converting that back to EMU (pts * 12700,
    # truncated) can land 1 EMU short of the original integer EMU value
    # because of the float32 precision loss. This finds the smallest
    # epsilon (in points) to add so the EMU conversion recovers the exact
    # original EMU value - i.e. it keeps text-only edits from nudging a
    # shape's stored position/size (which an autofit recalculation on text
    # assignment would otherwise do).
    $emuTarget = [int64](([double]$originalPts * 12700.0) + 0.5)
    for ($i = 0; $i -le 4000; $i++) {
        $eps = [Single]($i * 0.0000005)
        $candidate = [Single]($originalPts + $eps)
        $emu = [int64]([double]$candidate * 12700.0)
        if ($emu -ge $emuTarget) {
            return $eps
        }
    }
    return 0
}

function Save-ShapeGeometry($shp) {
    return @{
        Left   = $shp.Left
        Top    = $shp.Top
        Width  = $shp.Width
        Height = $shp.Height
    }
}

function Restore-ShapeGeometry($shp, $geo) {
    $shp.Left   = $geo.Left + (Get-RestoreEpsilon $geo.Left)
    $shp.Top    = $geo.Top + (Get-RestoreEpsilon $geo.Top)
    $shp.Width  = $geo.Width + (Get-RestoreEpsilon $geo.Width)
    $shp.Height = $geo.Height + (Get-RestoreEpsilon $geo.Height)
}

$p = $ppt.ActivePresentation

# --- Slide 4: "Fig 2. System Diagram" -> "Fig 3. System Diagram" -----------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(3)
$geo4 = Save-ShapeGeometry $shp4
$shp4.TextFrame.TextRange.Text = "Fig 3. System Diagram"
Restore-ShapeGeometry $shp4 $geo4

# --- Slide 5: "Fig 3. Tools (GCP, GitHub, Spark)" -> "Fig 4. Tools (GCP, GitHub, Spark)" ---
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(7)
$geo5 = Save-ShapeGeometry $shp5
$shp5.TextFrame.TextRange.Text = "Fig 4. Tools (GCP, GitHub, Spark)"
Restore-ShapeGeometry $shp5 $geo5

# --- Slide 6: "Fig 4. Beat Histogram for Classical (Left) and Pop (Right) [1]"
#     -> two runs: "Fig 5. " + "Beat Histogram for Classical (Left) and Pop (Right) [1]"
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(5)
$geo6 = Save-ShapeGeometry $shp6
$firstRun = $shp6.TextFrame.TextRange.Characters(1, 7)
$firstRun.Text = "Fig 5. "
Restore-ShapeGeometry $shp6 $geo6
